# Update purchase prices ("Kaufpreis $", column C) for two crypto assets:
#   row 4 -> SOL  : 37.38 -> 40
#   row 6 -> ATOM : 7     -> 7.5
# and leave the active selection on the last-edited cell (C6), matching the
# saved workbook view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = 40
$ws.Range("C6").Value = 7.5

$ws.Range("C6").Select()

$wb.Save()
